# Remove TDS from main invoice template
#
# The "Tax Deducted @ {meta:tds_tax_rate}" / "{meta:tds}" line lived in
# row 27 of Sheet1, between the SUB-TOTAL row (26) and the GRAND TOTAL
# row (previously 28). Deleting the entire row shifts every row below it
# up by one (GRAND TOTAL -> 27, blank spacer -> 28, "Total Amount (in
# words)" merge -> 29/30, signature note -> 31) and Excel automatically
# drops the now-unused shared strings for the TDS label/placeholder.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(27).Delete()
